$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 5 (old rows 5-10 shift down to 7-12)
$ws.Rows("5:6").Insert()

# New row 5: Mind_Muck (Demon, BattlecryMechanic)
$ws.Range("A5").Value = "Mind_Muck"
$ws.Range("B5").Value = 3
$ws.Range("C5").Value = 2
$ws.Range("D5").Value = "Minion"
$ws.Range("E5").Value = "Demon"
$ws.Range("F5").Value = 2
$ws.Range("G5").Value = 17
$ws.Range("H5").Value = "BattlecryMechanic"
$ws.Range("I5").Value = 1

# New row 6: Picky_Eater (Demon, BattlecryMechanic)
$ws.Range("A6").Value = "Picky_Eater"
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = 1
$ws.Range("D6").Value = "Minion"
$ws.Range("E6").Value = "Demon"
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 17
$ws.Range("H6").Value = "BattlecryMechanic"
$ws.Range("I6").Value = 0

$ws.Rows("5:6").RowHeight = 17

# The old Saltscale_Honcho row (now row 12) had its use_flg reset from 1 to 0
$ws.Range("I12").Value = 0

# Update selection / active cell to row 5 (whole row selected)
[void]$ws.Range("A5:XFD5").Select()
